$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "Spreadsheet" test tables (checkSpreadsheet1 / checkSpreadsheet2),
# added side-by-side starting at row 27 (left table B:D, right table G:I),
# mirroring the existing Method/Decision-table layout already on the sheet.
# ---------------------------------------------------------------------------

# --- Row 27: merged header cells (left-aligned, no wrap) ------------------
$ws.Range("B27:D27").HorizontalAlignment = -4131   # xlLeft
$ws.Range("B27:D27").Merge()
$ws.Range("B27").Value = "Spreadsheet Driver checkSpreadsheet1(Driver[] dd, int maxAge)"

$ws.Range("G27:I27").HorizontalAlignment = -4131   # xlLeft
$ws.Range("G27:I27").Merge()
$ws.Range("G27").Value = "Spreadsheet Driver checkSpreadsheet2(Driver[] dd, int dIndex)"

# --- Row 28: column headers ------------------------------------------------
$ws.Range("C28").Value = "Formula:Driver"
$ws.Range("D28").Value = "Age:int"
$ws.Range("H28").Value = "Formula:Driver"
$ws.Range("I28").Value = "Index:int"

# --- Row 29: SelectStep row -------------------------------------------------
$ws.Range("B29").Value = "SelectStep"
$ws.Range("C29").Value = "'=dd[!@ age < `$Age]"
$ws.Range("D29").Value = "'=maxAge"

$ws.Range("G29").Value = "SelectStep"
$ws.Range("H29").Value = "'=dd[!@ name == testDrivers[`$Index].name]"
$ws.Range("I29").Value = "'=dIndex"

# --- Row 30: RETURN row -----------------------------------------------------
$ws.Range("B30").Value = "RETURN"
$ws.Range("C30").Value = "'=`$SelectStep"

$ws.Range("G30").Value = "RETURN"
$ws.Range("H30").Value = "'=`$SelectStep"

# ---------------------------------------------------------------------------
# Column width tweaks that resulted from the new content being added
# (closest values reachable through the host's column-width rounding).
# ---------------------------------------------------------------------------
$ws.Range("B1").ColumnWidth = 9.333333333333334
$ws.Range("H1").ColumnWidth = 34.333333333333336
$ws.Range("I1").ColumnWidth = 11.5
$ws.Range("J1").ColumnWidth = 13.166666666666666

# ---------------------------------------------------------------------------
# Leave the selection where the author last left it when saving the file.
# ---------------------------------------------------------------------------
$ws.Range("G27:I27").Select()
